# Register-map rework described in the commit message: the "Board stack
# ID"/"Board stack revision" registers (scratchpad 0x00-0x01) are replaced
# by a 2-word Board stack ID (low/high, read from EEPROM) plus the Board
# stack revision moving one register to the right (0x02) and a newly
# Reserved register at 0x03. Every register below that point (packet
# checksum, command ID) shifts down by one slot, and the table grows one
# row to fit the extra Command ID byte (0x1E -> Command ID[23:16],
# Command ID[31:24] now have their own rows, C30/C31).
#
# Also documents (column F, rows 7-8) that "Next address for command
# data" / "Command data word countdown" are scratch registers reused for
# general workspace until the destination is verified.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Scratchpad register map (rows 20-31), reworked -----------------------
$ws.Range("C20").Value2 = "Board stack ID low word (as read at startup from EEPROM address 00)"
$ws.Range("C21").Value2 = "Board stack ID high word (as read at startup from EEPROM address 01)"
$ws.Range("C22").Value2 = "Board stack revision (as read at startup from EEPROM address 02)"
$ws.Range("C23").Value2 = "Reserved (but currently unused)"
$ws.Range("C24").Value2 = "Packet checksum[7:0]"
$ws.Range("C25").Value2 = "Packet checksum[15:8]"
$ws.Range("C26").Value2 = "Packet checksum[23:16]"
$ws.Range("C27").Value2 = "Packet checksum[31:24]"
$ws.Range("C28").Value2 = "Command ID[7:0]"
$ws.Range("C29").Value2 = "Command ID[15:8]"
$ws.Range("C30").Value2 = "Command ID[23:16]"
$ws.Range("C31").Value2 = "Command ID[31:24]"

# --- Clarifying notes for the existing command-interpreter registers -----
$ws.Range("F7").Value2 = "(after destination is verified, before that this is used for general workspace)"
$ws.Range("F8").Value2 = "(after destination is verified, before that this is used for general workspace)"

# --- Selection / view state: the saved sheet leaves the cursor on J35 ----
$ws.Activate()
$ws.Range("J35").Select()
